# Project Requirements list (numId=3) restructuring:
#  1. Merge the two runs of the first bullet ("Research and find ... applied to it")
#     into a single run.
#  2. Insert two new bullets after it ("Include an introduction section..." and
#     "Provide a brief list/description...").
#  3. Leave the next four existing bullets (about process summary, pandas,
#     sub-process functions, visualization) untouched content-wise - they just
#     shift down in the list.
#  4. Leave the "normally distributed" and "central tendency" bullets untouched.
#  5. Append two new bonus bullets at the end of the list.

$d = $word.ActiveDocument

$pkgHeader = '<?xml version="1.0" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">'
$pkgFooter = '</w:document></pkg:xmlData></pkg:part></pkg:package>'

$listPPr = '<w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="3"/></w:numPr><w:rPr><w:b/><w:bCs/><w:sz w:val="28"/><w:szCs w:val="28"/><w:u w:val="single"/></w:rPr></w:pPr>'

function Insert-ListParagraphAfter($paragraph, [string]$innerXml) {
    # Inserts a brand-new ListParagraph (numId=3) bullet immediately after
    # $paragraph, with body content $innerXml (the run/proofErr sequence).
    # Returns the newly created Word.Paragraph object.
    $paragraph.Range.InsertParagraphAfter()
    $newIndex = $paragraph.Index + 1
    $newPara = $d.Paragraphs.Item($newIndex)
    $start = $newPara.Range.Start
    $end = $newPara.Range.End
    $rng = $d.Range($start, $end)
    $xml = $pkgHeader + '<w:body><w:p>' + $listPPr + $innerXml + '</w:p><w:p/></w:body>' + $pkgFooter
    $rng.InsertXML($xml)
    return $d.Paragraphs.Item($newIndex)
}

# --- 1. Merge "Research and find a dataset that needs " + "cleaning processes
#        applied to it" into a single run -----------------------------------
$item1 = $d.Paragraphs.Item(12)
$s = $item1.Range.Start
$e = $item1.Range.End - 1
$rng = $d.Range($s, $e)
$rng.Text = "PLACEHOLDER_MERGE_TOKEN"
$item1b = $d.Paragraphs.Item(12)
$s2 = $item1b.Range.Start
$e2 = $item1b.Range.End - 1
$rng2 = $d.Range($s2, $e2)
$rng2.Text = "Research and find a dataset that needs cleaning processes applied to it"

# --- 2. Insert the two new requirement bullets after it ---------------------
$item1Final = $d.Paragraphs.Item(12)

$item2Inner = '<w:r><w:rPr><w:sz w:val="28"/><w:szCs w:val="28"/></w:rPr><w:t xml:space="preserve">Include an introduction section that gives a verbal explanation of the data set, where you found it (site the source), what </w:t></w:r><w:r><w:rPr><w:sz w:val="28"/><w:szCs w:val="28"/></w:rPr><w:t>the data describes/the purpose of the data, what you could potentially use it for, and a brief overview of the contents/purpose of your process/notebook</w:t></w:r>'
$item2Para = Insert-ListParagraphAfter $item1Final $item2Inner

$item3Inner = '<w:r><w:rPr><w:sz w:val="28"/><w:szCs w:val="28"/></w:rPr><w:t xml:space="preserve">Provide a brief list/description of the different skills/tools/methods that you learned in this program and applied in the notebook (make this list </w:t></w:r><w:proofErr w:type="gramStart"/><w:r><w:rPr><w:sz w:val="28"/><w:szCs w:val="28"/></w:rPr><w:t>as long as</w:t></w:r><w:proofErr w:type="gramEnd"/><w:r><w:rPr><w:sz w:val="28"/><w:szCs w:val="28"/></w:rPr><w:t xml:space="preserve"> you want but make it concise)</w:t></w:r>'
$item3Para = Insert-ListParagraphAfter $item2Para $item3Inner

# --- 3. Locate the "Calculate the 3 different measures..." bullet (it has not
#        moved relative to its neighbours, only the whole block shifted down
#        by 2 because of the two new bullets above) and append the two new
#        bonus bullets after it. ---------------------------------------------
$calcPara = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    if ($p.Range.Text -like "Calculate the 3 different measures*") {
        $calcPara = $p
        break
    }
}

$item10Inner = '<w:r><w:rPr><w:sz w:val="28"/><w:szCs w:val="28"/></w:rPr><w:lastRenderedPageBreak/><w:t xml:space="preserve">Bonus: implement some or </w:t></w:r><w:proofErr w:type="gramStart"/><w:r><w:rPr><w:sz w:val="28"/><w:szCs w:val="28"/></w:rPr><w:t>all of</w:t></w:r><w:proofErr w:type="gramEnd"/><w:r><w:rPr><w:sz w:val="28"/><w:szCs w:val="28"/></w:rPr><w:t xml:space="preserve"> your functions via organized modules that you import from different file(s) (organized meaning all stats function in one file and all data cleaning processes in another)</w:t></w:r>'
$item10Para = Insert-ListParagraphAfter $calcPara $item10Inner

$item11Inner = '<w:r><w:rPr><w:sz w:val="28"/><w:szCs w:val="28"/></w:rPr><w:t xml:space="preserve">Bonus: provide a summary of how you think you could use </w:t></w:r><w:r><w:rPr><w:sz w:val="28"/><w:szCs w:val="28"/></w:rPr><w:t xml:space="preserve">the data to create predictions, calculations, or whatever </w:t></w:r><w:proofErr w:type="gramStart"/><w:r><w:rPr><w:sz w:val="28"/><w:szCs w:val="28"/></w:rPr><w:t>future plans</w:t></w:r><w:proofErr w:type="gramEnd"/><w:r><w:rPr><w:sz w:val="28"/><w:szCs w:val="28"/></w:rPr><w:t xml:space="preserve"> you potentially think you could implement (these don&#8217;t need to be too detailed for now; this is just to get you thinking &#8211; you can add on later)</w:t></w:r>'
$item11Para = Insert-ListParagraphAfter $item10Para $item11Inner

Write-Host "Done. Paragraph count now: $($d.Paragraphs.Count)"
